$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B" + $r1 + ":AB" + $r1)
    $rng2 = $ws.Range("B" + $r2 + ":AB" + $r2)

    $v1 = $rng1.Value()
    $v2 = $rng2.Value()

    $rng1.Value = $v2
    $rng2.Value = $v1
}

Swap-Rows 29 30
Swap-Rows 147 148
Swap-Rows 173 174
